$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "C"="0.04864846792477806"; "D"="0.1941765291412167"; "E"="0.1678980535918484"; "F"="1.420009619434879"; "G"="0.002460955856352731"; "J"="0.1875123552826565"; "K"="1.984762760414014"; "O"="3.418384124897045" }
  3 = @{ "C"="0.04322810012756406"; "D"="0.1874186996394798"; "E"="0.1634744811215398"; "F"="1.429637574686993"; "G"="0.002464427855176485"; "J"="0.1836935370985913"; "K"="1.771614680571133"; "O"="3.460778934695327" }
  4 = @{ "C"="0.03991531886543953"; "D"="0.1833153322243959"; "E"="0.1608344221120284"; "F"="1.436672043547993"; "G"="0.002466670940765399"; "J"="0.1814652300583433"; "K"="1.640514740126662"; "O"="3.489828238116857" }
  5 = @{ "C"="0.0385691550447973"; "D"="0.1816548435506604"; "E"="0.1597777205352919"; "F"="1.439820357181794"; "G"="0.002467613081776093"; "J"="0.1805864147440417"; "K"="1.58703640985658"; "O"="3.502422980041047" }
  6 = @{ "C"="0.03834585528869638"; "D"="0.1813798285321297"; "E"="0.1596034130449979"; "F"="1.440360128053698"; "G"="0.002467771221261854"; "J"="0.1804422526538971"; "K"="1.578153194308584"; "O"="3.504559979669509" }
  7 = @{ "C"="0.03989714860524884"; "D"="0.1832928909092431"; "E"="0.1608200935177351"; "F"="1.436713363030918"; "G"="0.002466683533047262"; "J"="0.1814532597265242"; "K"="1.639793727678693"; "O"="3.489995033384176" }
  8 = @{ "C"="0.04677631561888518"; "D"="0.1918369510690781"; "E"="0.1663570382808857"; "F"="1.423095898802856"; "G"="0.0024621299610361"; "J"="0.1861714244319401"; "K"="1.911317858290829"; "O"="3.432374017467197" }
  9 = @{ "C"="0.060390502675844"; "D"="0.2089528160980052"; "E"="0.1778179313800763"; "F"="1.405331610978365"; "G"="0.002454079269666277"; "J"="0.1963507164034723"; "K"="2.441885180827967"; "O"="3.343429728207326" }
  10 = @{ "C"="0.07047286026893573"; "D"="0.221744368058836"; "E"="0.1866065463366908"; "F"="1.397773163144322"; "G"="0.002448694592733589"; "J"="0.2043996737552618"; "K"="2.830452470996534"; "O"="3.292880653044818" }
  11 = @{ "C"="0.07507795240437076"; "D"="0.2276099576311026"; "E"="0.1906849254740948"; "F"="1.395536798181269"; "G"="0.002446358888686087"; "J"="0.2081863509418156"; "K"="3.006937157757989"; "O"="3.273127375713244" }
  12 = @{ "C"="0.07682451019792325"; "D"="0.229837729524462"; "E"="0.1922408587982218"; "F"="1.394863536550631"; "G"="0.00244549069394034"; "J"="0.2096383395635684"; "K"="3.073725471534487"; "O"="3.266115897412021" }
  13 = @{ "C"="0.07644823662143096"; "D"="0.2293576471275429"; "E"="0.1919052478131036"; "F"="1.395000802515924"; "G"="0.002445676952057719"; "J"="0.2093248238613938"; "K"="3.05934335280881"; "O"="3.267605061145531" }
  14 = @{ "C"="0.07522158824734504"; "D"="0.2277931060936425"; "E"="0.1908127018900601"; "F"="1.395477924994722"; "G"="0.002446287135688598"; "J"="0.2083054447163164"; "K"="3.012432742525334"; "O"="3.272541130222663" }
  15 = @{ "C"="0.07447058466756573"; "D"="0.2268356359389543"; "E"="0.1901449881626007"; "F"="1.395792806676212"; "G"="0.002446663009467333"; "J"="0.2076833987151048"; "K"="2.983693007086572"; "O"="3.275625722527536" }
  16 = @{ "C"="0.07017228372679085"; "D"="0.2213619674397478"; "E"="0.1863416305975036"; "F"="1.397943566779489"; "G"="0.002448849519841239"; "J"="0.2041547290519787"; "K"="2.818912973020417"; "O"="3.294237022078278" }
  17 = @{ "C"="0.06754020031053187"; "D"="0.2180159198544658"; "E"="0.1840289696465547"; "F"="1.399571388147294"; "G"="0.002450219965890788"; "J"="0.2020221032046692"; "K"="2.717752894856915"; "O"="3.306486519249518" }
  18 = @{ "C"="0.06602804182769262"; "D"="0.2160957588874197"; "E"="0.182706356513421"; "F"="1.400620761775002"; "G"="0.002451018927557173"; "J"="0.2008072509830896"; "K"="2.659542432368596"; "O"="3.31383708532411" }
  19 = @{ "C"="0.06551634971974352"; "D"="0.2154463838860323"; "E"="0.1822598429407165"; "F"="1.400995462692549"; "G"="0.002451291285351487"; "J"="0.2003979439960801"; "K"="2.639829015822158"; "O"="3.316378165094079" }
  20 = @{ "C"="0.06782020907712649"; "D"="0.2183716580335755"; "E"="0.1842743732193881"; "F"="1.399386393534456"; "G"="0.002450072970786588"; "J"="0.2022479055002719"; "K"="2.728524255174989"; "O"="3.305150957034044" }
  21 = @{ "C"="0.07558181077490644"; "D"="0.2282524714368321"; "E"="0.1911332960912233"; "F"="1.395333065058509"; "G"="0.002446107468673365"; "J"="0.208604370581881"; "K"="3.026212704622708"; "O"="3.271078547876641" }
  22 = @{ "C"="0.08067027611038213"; "D"="0.2347485776889187"; "E"="0.1956832515025297"; "F"="1.393696211939755"; "G"="0.002443610672571974"; "J"="0.2128639601493347"; "K"="3.220519688785885"; "O"="3.25154312463215" }
  23 = @{ "C"="0.07795300716981046"; "D"="0.2312780029040766"; "E"="0.19324870773508"; "F"="1.394476962784893"; "G"="0.002444934603767157"; "J"="0.2105808854793878"; "K"="3.116838160598547"; "O"="3.261718663818925" }
  24 = @{ "C"="0.06769361381330441"; "D"="0.2182108178124196"; "E"="0.1841634045645435"; "F"="1.39946967604628"; "G"="0.00245013939288371"; "J"="0.2021457853310835"; "K"="2.723654685526697"; "O"="3.305753804946676" }
  25 = @{ "C"="0.05669376618271826"; "D"="0.2042842470229118"; "E"="0.1746528340273841"; "F"="1.409175821444236"; "G"="0.002456163689288478"; "J"="0.1934971814717699"; "K"="2.298564372189901"; "O"="3.364902651233393" }
}

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  foreach ($col in $rowVals.Keys) {
    $ws.Range("$col$r").Value = [double]$rowVals[$col]
  }
}